$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new uploaded documents (Hangsihak_Sin_CoverLetter_2.pdf, Hangsihak_Sin_Transcript_1.pdf)
# are compared alongside the two already-present ones (net2.pdf, Test_Plan_Requirements.pdf),
# growing the similarity matrix from 2x2 to 4x4. The newly-uploaded docs take the first
# header column / row-label slot, pushing the previously-existing docs to columns D/E and
# rows 4/5.

$newDoc1 = "C:/Users/Hangsihak Sin/OneDrive/Documents/School/Doc-Wise/backend/phase_one/temp_files\Hangsihak_Sin_CoverLetter_2.pdf"
$newDoc2 = "C:/Users/Hangsihak Sin/OneDrive/Documents/School/Doc-Wise/backend/phase_one/temp_files\Hangsihak_Sin_Transcript_1.pdf"
$doc3 = "C:/Users/Hangsihak Sin/OneDrive/Documents/School/Doc-Wise/backend/phase_one/temp_files\net2.pdf"
$doc4 = "C:/Users/Hangsihak Sin/OneDrive/Documents/School/Doc-Wise/backend/phase_one/temp_files\Test_Plan_Requirements.pdf"

# New header cells / row labels for the previously-existing docs, shifted out to D1:E1 / A4:A5
$ws.Range("D1").Value = $doc3
$ws.Range("E1").Value = $doc4
$ws.Range("A4").Value = $doc3
$ws.Range("A5").Value = $doc4

# Give the new cells the same bold/centered/thin-border look as the existing header cells,
# copied in one shot so no stray intermediate formats are minted.
$ws.Range("B1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)

# The first header column / row label now hold the newly-uploaded documents
$ws.Range("B1").Value = $newDoc1
$ws.Range("C1").Value = $newDoc2
$ws.Range("A2").Value = $newDoc1
$ws.Range("A3").Value = $newDoc2

# Updated similarity matrix (4x4)
$ws.Range("B2").Value = 0.9999999403953553
$ws.Range("C2").Value = 0.4485565959896671
$ws.Range("D2").Value = 0.4427806619818045
$ws.Range("E2").Value = 0.431584213614177

$ws.Range("B3").Value = 0.4485565959896671
$ws.Range("C3").Value = 0.9999998807907116
$ws.Range("D3").Value = 0.392725977505771
$ws.Range("E3").Value = 0.3869618633737347

$ws.Range("B4").Value = 0.4427806619818045
$ws.Range("C4").Value = 0.392725977505771
$ws.Range("D4").Value = 0.9999999403953526
$ws.Range("E4").Value = 0.4787833090451574

$ws.Range("B5").Value = 0.431584213614177
$ws.Range("C5").Value = 0.3869618633737347
$ws.Range("D5").Value = 0.4787833090451574
$ws.Range("E5").Value = 1
